# Auto-generated edit script applying the crypto price/volume update diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to keep the assigned string as text (avoid Excel
    # auto-converting numeric-looking strings like "0.998" or "3.10" into
    # real numbers, which would silently drop formatting / trailing zeros),
    # then restore the original (default) cell style so no stray number
    # format is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.652.23"
Set-TextValue $ws.Range("E2") "  +0.36%  "
Set-TextValue $ws.Range("D3") "3.496.41"
Set-TextValue $ws.Range("E3") "  +0.24%  "
Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "603.61"
Set-TextValue $ws.Range("E5") "  -0.98%  "
Set-TextValue $ws.Range("D6") "194.31"
Set-TextValue $ws.Range("E6") "  +4.48%  "
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("E8") "  -0.06%  "
Set-TextValue $ws.Range("E9") "  -6.48%  "
Set-TextValue $ws.Range("D10") "0.648"
Set-TextValue $ws.Range("E10") "  +0.18%  "
Set-TextValue $ws.Range("D11") "53.21"
Set-TextValue $ws.Range("E11") "  +0.40%  "
Set-TextValue $ws.Range("D13") "9.45"
Set-TextValue $ws.Range("E13") "  -0.75%  "
Set-TextValue $ws.Range("D14") "4.063.53"
Set-TextValue $ws.Range("E14") "  +0.49%  "
Set-TextValue $ws.Range("D15") "593.38"
Set-TextValue $ws.Range("E15") "  -1.31%  "
Set-TextValue $ws.Range("D16") "69.810.06"
Set-TextValue $ws.Range("E16") "  +0.46%  "
Set-TextValue $ws.Range("E17") "  +1.45%  "
Set-TextValue $ws.Range("D18") "18.97"
Set-TextValue $ws.Range("E18") "  +0.66%  "
Set-TextValue $ws.Range("B19") "TRON"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D19") "0.122"
Set-TextValue $ws.Range("E19") "  +1.45%  "
Set-TextValue $ws.Range("B20") "WrappedEther"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D20") "3.489.50"
Set-TextValue $ws.Range("E20") "  -0.55%  "
Set-TextValue $ws.Range("D21") "0.984"
Set-TextValue $ws.Range("E21") "  -0.07%  "
Set-TextValue $ws.Range("D22") "18.04"
Set-TextValue $ws.Range("E22") "  +5.23%  "
Set-TextValue $ws.Range("D23") "5.32"
Set-TextValue $ws.Range("E23") "  +4.21%  "
Set-TextValue $ws.Range("D24") "4.65"
Set-TextValue $ws.Range("E24") "  +0.52%  "
Set-TextValue $ws.Range("D25") "101.93"
Set-TextValue $ws.Range("E25") "  -3.79%  "
Set-TextValue $ws.Range("E26") "  +2.98%  "
Set-TextValue $ws.Range("D27") "10.84"
Set-TextValue $ws.Range("E27") "  -1.00%  "
Set-TextValue $ws.Range("E28") "  -1.63%  "
Set-TextValue $ws.Range("E29") "  -0.75%  "
Set-TextValue $ws.Range("D30") "7.01"
Set-TextValue $ws.Range("E30") "  +0.54%  "
Set-TextValue $ws.Range("D31") "4.21"
Set-TextValue $ws.Range("E31") "  +1.69%  "
Set-TextValue $ws.Range("D32") "12.38"
Set-TextValue $ws.Range("E32") "  -0.37%  "
Set-TextValue $ws.Range("D33") "0.115"
Set-TextValue $ws.Range("E33") "  -0.44%  "
Set-TextValue $ws.Range("E34") "  -0.41%  "
Set-TextValue $ws.Range("D35") "0.0₃0829"
Set-TextValue $ws.Range("E35") "  +6.57%  "
Set-TextValue $ws.Range("D36") "3.711.01"
Set-TextValue $ws.Range("E36") "  +3.16%  "
Set-TextValue $ws.Range("D37") "3.10"
Set-TextValue $ws.Range("E37") "  -2.40%  "
Set-TextValue $ws.Range("D38") "0.999"
Set-TextValue $ws.Range("E38") "  +0.10%  "
Set-TextValue $ws.Range("E39") "  -1.08%  "
Set-TextValue $ws.Range("D40") "0.390"
Set-TextValue $ws.Range("E40") "  -1.56%  "
Set-TextValue $ws.Range("D41") "36.33"
Set-TextValue $ws.Range("E41") "  -1.01%  "
Set-TextValue $ws.Range("D42") "478.72"
Set-TextValue $ws.Range("E42") "  -8.21%  "
Set-TextValue $ws.Range("E43") "  -2.72%  "
Set-TextValue $ws.Range("D44") "0.0452"
Set-TextValue $ws.Range("E44") "  -2.01%  "
Set-TextValue $ws.Range("E45") "  -1.85%  "
Set-TextValue $ws.Range("E46") "  -4.67%  "
Set-TextValue $ws.Range("D47") "3.28"
Set-TextValue $ws.Range("E47") "  -1.38%  "
Set-TextValue $ws.Range("E48") "  +0.23%  "
Set-TextValue $ws.Range("D49") "8.41"
Set-TextValue $ws.Range("E49") "  -4.38%  "
Set-TextValue $ws.Range("D50") "0.000244"
Set-TextValue $ws.Range("E50") "  +1.01%  "
Set-TextValue $ws.Range("E51") "  +9.96%  "
